# Atualiza pontuacoes e resultados das competicoes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 64.5
$ws.Range("B3").Value  = 63.56
$ws.Range("B4").Value  = 47.86
$ws.Range("B5").Value  = 66.37
$ws.Range("B7").Value  = 53.66
$ws.Range("B8").Value  = 71.95999999999999
$ws.Range("B9").Value  = 56.05
$ws.Range("B10").Value = 61.56
$ws.Range("B12").Value = 47.86
$ws.Range("B13").Value = 55.66
$ws.Range("B14").Value = 62.56
$ws.Range("B16").Value = 59.8
$ws.Range("B17").Value = 61.16
$ws.Range("B19").Value = 54.9
$ws.Range("B20").Value = 68.06
